# Natmi following Dr Hou advice
# Recomputed NATMI ligand-receptor edge stats for Calr-Tshr: the target
# cluster set was expanded so every sending cluster (ECs, FAPs, M2, sCs)
# now has a row for every target cluster (ECs, FAPs, M2, sCs) -> 16 data
# rows (2-17) instead of 12 (2-13), and all numeric columns (E:T) were
# refreshed with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @("ECs", "Calr", "Tshr", "ECs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 58.42598751752668, 525.8338876577401, 0.02503815794992931, 0.02503815794992931)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $rowData[$i] }

$rowData = @("ECs", "Calr", "Tshr", "FAPs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 274.4947006480564, 2470.452305832508, 0.1176332992092411, 0.1176332992092411)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $rowData[$i] }

$rowData = @("ECs", "Calr", "Tshr", "M2", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 9.460965994219112, 85.14869394797202, 0.00405444855940359, 0.004054448559403589)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $rowData[$i] }

$rowData = @("ECs", "Calr", "Tshr", "sCs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 38.65933269313112, 347.93399423818, 0.01656725918272437, 0.01656725918272436)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $rowData[$i] }

$rowData = @("FAPs", "Calr", "Tshr", "ECs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 128.59907640679, 1157.39168766111, 0.05511047607577593, 0.05511047607577592)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $rowData[$i] }

$rowData = @("FAPs", "Calr", "Tshr", "FAPs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 604.1791757701847, 5437.612581931662, 0.2589178938302757, 0.2589178938302757)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $rowData[$i] }

$rowData = @("FAPs", "Calr", "Tshr", "M2", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 20.82414932922867, 187.417343963058, 0.008924082625419594, 0.008924082625419594)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $rowData[$i] }

$rowData = @("FAPs", "Calr", "Tshr", "sCs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 85.09149250319668, 765.8234325287701, 0.03646552364820783, 0.03646552364820783)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $rowData[$i] }

$rowData = @("M2", "Calr", "Tshr", "ECs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 105.4298747153783, 948.8688724384052, 0.04518143326157768, 0.04518143326157767)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $rowData[$i] }

$rowData = @("M2", "Calr", "Tshr", "FAPs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 495.3265341159779, 4457.938807043802, 0.2122696513134078, 0.2122696513134078)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $rowData[$i] }

$rowData = @("M2", "Calr", "Tshr", "M2", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 17.07234232297322, 153.651080906759, 0.007316264933128206, 0.007316264933128205)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $rowData[$i] }

$rowData = @("M2", "Calr", "Tshr", "sCs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 69.76088510603724, 627.8479659543351, 0.02989567030403825, 0.02989567030403825)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $rowData[$i] }

$rowData = @("sCs", "Calr", "Tshr", "ECs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 65.34318201552001, 588.08863813968, 0.02800248625262483, 0.02800248625262482)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $rowData[$i] }

$rowData = @("sCs", "Calr", "Tshr", "FAPs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 306.992794625184, 2762.935151626656, 0.131560191070964, 0.131560191070964)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $rowData[$i] }

$rowData = @("sCs", "Calr", "Tshr", "M2", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 10.581072725856, 95.22965453270402, 0.004534464566990803, 0.004534464566990802)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $rowData[$i] }

$rowData = @("sCs", "Calr", "Tshr", "sCs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 43.23630494064, 389.1267444657601, 0.01852869721629108, 0.01852869721629107)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(17, $i + 1).Value = $rowData[$i] }
